# Applies the cryptos list refresh described in the commit:
# "Updated cryptos list on Fri Mar  3 13:49:03 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '22.342.54'
$ws.Cells.Item(2, 5).Value = '  -4.54%  '

$ws.Cells.Item(3, 4).Value = '1.567.79'
$ws.Cells.Item(3, 5).Value = '  -4.19%  '

$ws.Cells.Item(4, 4).Value = '''1.002'
$ws.Cells.Item(4, 5).Value = '  +0.07%  '

$ws.Cells.Item(5, 4).Value = '''1.002'
$ws.Cells.Item(5, 5).Value = '  +0.07%  '

$ws.Cells.Item(6, 4).Value = '''289.08'
$ws.Cells.Item(6, 5).Value = '  -3.63%  '

$ws.Cells.Item(7, 4).Value = '''0.3689'
$ws.Cells.Item(7, 5).Value = '  -2.41%  '

$ws.Cells.Item(8, 4).Value = '''49.28'
$ws.Cells.Item(8, 5).Value = '  -2.00%  '

$ws.Cells.Item(9, 4).Value = '''0.3373'
$ws.Cells.Item(9, 5).Value = '  -4.30%  '

$ws.Cells.Item(10, 4).Value = '''1.161'
$ws.Cells.Item(10, 5).Value = '  -3.94%  '

$ws.Cells.Item(11, 4).Value = '''0.07592'
$ws.Cells.Item(11, 5).Value = '  -5.76%  '

$ws.Cells.Item(12, 5).Value = '  +0.09%  '

$ws.Cells.Item(13, 4).Value = '''21.14'
$ws.Cells.Item(13, 5).Value = '  -3.85%  '

$ws.Cells.Item(14, 4).Value = '''6.042'
$ws.Cells.Item(14, 5).Value = '  -4.67%  '

$ws.Cells.Item(15, 4).Value = '''6.882'
$ws.Cells.Item(15, 5).Value = '  -5.17%  '

$ws.Cells.Item(16, 4).Value = '1.566.98'
$ws.Cells.Item(16, 5).Value = '  -4.47%  '

$ws.Cells.Item(17, 4).Value = '''0.00001131'
$ws.Cells.Item(17, 5).Value = '  -6.02%  '

$ws.Cells.Item(18, 4).Value = '''88.88'
$ws.Cells.Item(18, 5).Value = '  -7.40%  '

$ws.Cells.Item(19, 5).Value = '  -2.82%  '

$ws.Cells.Item(20, 5).Value = '  -0.02%  '

$ws.Cells.Item(21, 4).Value = '''6.228'
$ws.Cells.Item(21, 5).Value = '  -6.85%  '

$ws.Cells.Item(22, 2).Value = 'BitDAO'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit'
$ws.Cells.Item(22, 4).Value = '''0.5321'
$ws.Cells.Item(22, 5).Value = '  -7.35%  '

$ws.Cells.Item(23, 2).Value = 'Avalanche'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Cells.Item(23, 4).Value = '''16.52'
$ws.Cells.Item(23, 5).Value = '  -4.69%  '

$ws.Cells.Item(24, 4).Value = '''11.92'
$ws.Cells.Item(24, 5).Value = '  -3.18%  '

$ws.Cells.Item(25, 4).Value = '22.366.59'
$ws.Cells.Item(25, 5).Value = '  -4.48%  '

$ws.Cells.Item(26, 4).Value = '''2.384'
$ws.Cells.Item(26, 5).Value = '  -3.37%  '

$ws.Cells.Item(27, 4).Value = '''2.973'
$ws.Cells.Item(27, 5).Value = '  +2.67%  '

$ws.Cells.Item(28, 4).Value = '''19.93'
$ws.Cells.Item(28, 5).Value = '  -4.26%  '

$ws.Cells.Item(29, 4).Value = '''145.18'
$ws.Cells.Item(29, 5).Value = '  -4.43%  '

$ws.Cells.Item(30, 4).Value = '''4.962'
$ws.Cells.Item(30, 5).Value = '  -4.36%  '

$ws.Cells.Item(31, 4).Value = '''125.18'
$ws.Cells.Item(31, 5).Value = '  -5.64%  '

$ws.Cells.Item(32, 4).Value = '1.745.42'
$ws.Cells.Item(32, 5).Value = '  -4.12%  '

$ws.Cells.Item(33, 4).Value = '''1.036'
$ws.Cells.Item(33, 5).Value = '  +6.14%  '

$ws.Cells.Item(34, 4).Value = '''6.234'
$ws.Cells.Item(34, 5).Value = '  -8.84%  '

$ws.Cells.Item(35, 4).Value = '''1.993'
$ws.Cells.Item(35, 5).Value = '  -6.27%  '

$ws.Cells.Item(36, 4).Value = '''10.27'
$ws.Cells.Item(36, 5).Value = '  -9.40%  '

$ws.Cells.Item(37, 4).Value = '''0.08436'
$ws.Cells.Item(37, 5).Value = '  -3.24%  '

$ws.Cells.Item(38, 4).Value = '''0.02526'

$ws.Cells.Item(39, 4).Value = '''0.2326'
$ws.Cells.Item(39, 5).Value = '  -3.96%  '

$ws.Cells.Item(40, 4).Value = '''5.526'
$ws.Cells.Item(40, 5).Value = '  -6.04%  '

$ws.Cells.Item(41, 4).Value = '''0.06495'
$ws.Cells.Item(41, 5).Value = '  -4.42%  '

$ws.Cells.Item(42, 4).Value = '''11.75'
$ws.Cells.Item(42, 5).Value = '  -9.45%  '

$ws.Cells.Item(43, 4).Value = '''1.240'
$ws.Cells.Item(43, 5).Value = '  -4.59%  '

$ws.Cells.Item(44, 4).Value = '''0.6347'
$ws.Cells.Item(44, 5).Value = '  -7.21%  '

$ws.Cells.Item(45, 4).Value = '''14.29'
$ws.Cells.Item(45, 5).Value = '  -8.50%  '

$ws.Cells.Item(46, 4).Value = '''1.001'
$ws.Cells.Item(46, 5).Value = '  +0.01%  '

$ws.Cells.Item(47, 4).Value = '''0.5965'
$ws.Cells.Item(47, 5).Value = '  -5.62%  '

$ws.Cells.Item(48, 4).Value = '''3.749'
$ws.Cells.Item(48, 5).Value = '  -3.95%  '

$ws.Cells.Item(49, 4).Value = '''2.120'
$ws.Cells.Item(49, 5).Value = '  -5.35%  '

$ws.Cells.Item(50, 4).Value = '''1.247'
$ws.Cells.Item(50, 5).Value = '  +2.59%  '

$ws.Cells.Item(51, 4).Value = '''123.02'
$ws.Cells.Item(51, 5).Value = '  -3.07%  '
